$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Extend the "blank filler" formatting (the plain bordered/filled
#    style already used by B39) down across the newly appended blank
#    rows 41-50, columns B:G.
# ------------------------------------------------------------------
$null = $ws.Range("B39").Copy()
$null = $ws.Range("B41:G50").PasteSpecial(-4122)   # xlPasteFormats

# ------------------------------------------------------------------
# 2. Give rows 39 and 40 the same per-column formatting pattern used
#    by the existing data rows (row 4 is a normal "expense" row).
# ------------------------------------------------------------------
$null = $ws.Range("B4:G4").Copy()
$null = $ws.Range("B39:G40").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Fill in the new transaction data.
# ------------------------------------------------------------------
# Row 39: 2018-05-17, 支出, 其他, 班费100+拿药100, 200
$ws.Range("C39").Value = "支出"
$ws.Range("D39").Value = 200
$ws.Range("E39").Value = 43237
$ws.Range("F39").Value = "其他"
$ws.Range("G39").Value = "班费100+拿药100"

# Row 40: 2018-05-21, 支出, 生活费, 生活费(5/21-5/31), 400
$ws.Range("C40").Value = "支出"
$ws.Range("D40").Value = 400
$ws.Range("E40").Value = 43241
$ws.Range("F40").Value = "生活费"
$ws.Range("G40").Value = "生活费(5/21-5/31)"

# ------------------------------------------------------------------
# 4. Update the view: scroll down a bit and leave the final selection
#    on K36, matching where the user ended up after the edit.
# ------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$null = $ws.Range("K36").Select()
